$d = $word.ActiveDocument

# --- Edit 1: remove the _GoBack bookmark from its original location
# (right after "==0); {") ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Edit 2: replace the two empty paragraphs that follow the question
# "¿Cómo funciona el programa que muestra la figura de la derecha?"
# with the new answer text (and re-add the _GoBack bookmark on the
# trailing empty paragraph) ---
$rng = $d.Content
$found = $rng.Find.Execute("Cómo funciona el programa que muestra la figura de la derecha?", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate question 19 text in the document."
}

$rng.Collapse(0)
$questionPara = $rng.Paragraphs(1)
$emptyPara1 = $questionPara.Next()
$emptyPara2 = $emptyPara1.Next()

$newPara1Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:ind w:left="360"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cs="Arial"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">La condición cuando k=0 la ejecuta solo una vez porque cada que el </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cs="Arial"/></w:rPr><w:t>for</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cs="Arial"/></w:rPr><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cs="Arial"/></w:rPr><w:t>;;) se ejecuta k aumenta y a medida que esto pasa sale el mensaje que aparece en código y se hace la multiplicación de k*k.</w:t></w:r></w:p>'

$newPara2Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:ind w:left="360"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cs="Arial"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

[void]$emptyPara2.Range.InsertXML($newPara2Xml)
[void]$emptyPara1.Range.InsertXML($newPara1Xml)

Write-Host "Done"
